# Zeiten aktualisiert und Sound beim essen einer Beere
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Append the new log rows (31-40) using the OLD column layout
#    (A=Datum, B=Von, C=Bis, D=Mitglied, E=Taetigkeit) so the new shared
#    strings get interned in the same order as the source edit, *before*
#    the new "Dauer" column is inserted.
# ---------------------------------------------------------------------------
$ws.Range("A30:E30").Copy()
$ws.Range("A31:E40").PasteSpecial(-4122)

$newRows = @(
    @(42136, 0.33333333333333331, 0.46875,               "Juliano, Tommy, Constantin, Jot", "Schwierigkeitsanpassung im Singleplayer und Responsive Layout eingebaut"),
    @(42138, 0.66666666666666663, 0.875,                 "Tommy",                            "Threads zur Bluetoothverbindung eingebaut"),
    @(42145, 0.92708333333333337, 0.010416666666666666, "Tommy",                            "BluetoothService erstellt und GPMultiPlayer Activity erstellt"),
    @(42149, 0.04513888888888889, 0.11319444444444444,  "Tommy",                            "Bluetoothverbindung zwischen zwei Geräten geht jetzt"),
    @(42150, 0.95833333333333337, 0.97916666666666663,  "Tommy",                            "Unterscheidung zwischen ersten und zweiten Player ist jetzt möglich"),
    @(42155, 0.0763888888888889, 0.10555555555555556,   "Tommy",                            "Positionen des ersten Players, werden jetzt an den zweiten Player gesendet"),
    @(42163, 0.625, 0.67291666666666661,                 "Juliano",                          "Schlange wird in extra Thread der View gezeichnet"),
    @(42164, 0.75, 0.81597222222222221,                  "Tommy",                            "Anpassungen in MultiplayerView Framework"),
    @(42165, 0.33333333333333331, 0.99722222222222223,  "Juliano, Tommy, Constantin, Jot", "Beide Schlangen werden jetzt gezeichnet"),
    @(42166, 0.39583333333333331, 0.56527777777777777,  "Juliano, Tommy, Constantin, Jot", "Multiplayer beendet bei Kollision")
)

$r = 31
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Insert the new "Dauer" column at D (shifts old D/E/F -> E/F/G)
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = 8.83

# Header
$ws.Range("D1").Value = "Gesamt"

# ---------------------------------------------------------------------------
# 3) Fill the duration column. Most rows compute "=Cn-Bn"; a handful of rows
#    cross midnight (Bis < Von) so the original author typed the literal
#    duration instead of a formula.
# ---------------------------------------------------------------------------
$literalDauer = @{
    5  = 0.041666666666666664
    16 = 0.14027777777777778
    28 = 0.04861111111111111
    29 = 0.3125
    30 = 0.15277777777777776
    33 = 0.08333333333333333
}

for ($row = 3; $row -le 40; $row++) {
    if ($row -eq 2) { continue }
    if ($literalDauer.ContainsKey($row)) {
        $ws.Range("D$row").Value = $literalDauer[$row]
    } else {
        $ws.Range("D$row").Formula = "=C$row-B$row"
    }
}

# ---------------------------------------------------------------------------
# 4) Blank separator row + the four "Gesamt" summary rows
# ---------------------------------------------------------------------------
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").ClearContents()

$ws.Range("A42").Value = "Juliano"
$ws.Range("B42").Formula = "=D3+D7+D8+D9+D11+D15+D17+D19+D21+D22+D24+D25+D26+D27+D28+D29+D31+D37+D39+D40"

$ws.Range("A43").Value = "Tommy"
$ws.Range("B43").Formula = "=+D3+D4+D5+D10+D9+D11+D13+D14+D15+D16+D18+D19+D20+D22+D27+D30+D31+D32+D33+D34+D35+D36+D38+D39+D40"

$ws.Range("A44").Value = "Constantin"
$ws.Range("B44").Formula = "=D3+D6+D9+D11+D12+D15+D19+D22+D27+D31+D39+D40"

$ws.Range("A45").Value = "Jotprabh"
$ws.Range("B45").Formula = "=D3+D9+D11+D15+D19+D22+D27+D31+D39+D40"

# ---------------------------------------------------------------------------
# 5) View bookkeeping: restore the scroll position / selection the author
#    left the sheet in.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C45").Select()
